$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.987.54'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '2.244.31'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '97.29'
$ws.Range("E5").Value = '  +17.08%  '
$ws.Range("D6").Value = '272.15'
$ws.Range("E6").Value = '  +5.16%  '
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +5.84%  '
$ws.Range("D10").Value = '47.82'
$ws.Range("E10").Value = '  +8.01%  '
$ws.Range("D11").Value = '0.0941'
$ws.Range("E11").Value = '  +2.37%  '
$ws.Range("D12").Value = '8.32'
$ws.Range("E12").Value = '  +16.19%  '
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.30'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +6.67%  '
$ws.Range("D15").Value = '2.577.05'
$ws.Range("E15").Value = '  +2.01%  '
$ws.Range("D16").Value = '0.828'
$ws.Range("E16").Value = '  +5.74%  '
$ws.Range("D17").Value = '2.240.92'
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = '43.923.82'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("E19").Value = '  +2.28%  '
$ws.Range("E20").Value = '  +4.86%  '
$ws.Range("D21").Value = '70.87'
$ws.Range("D22").Value = '2.34'
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").Value = '234.48'
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("D24").Value = '9.39'
$ws.Range("E24").Value = '  +4.70%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '11.43'
$ws.Range("E26").Value = '  +7.42%  '
$ws.Range("E27").Value = '  +11.69%  '
$ws.Range("E28").Value = '  +2.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.60'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.23%  '
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("D31").Value = '173.45'
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  +6.00%  '
$ws.Range("D33").Value = '21.11'
$ws.Range("E33").Value = '  +3.61%  '
$ws.Range("D34").Value = '5.61'
$ws.Range("E34").Value = '  +5.62%  '
$ws.Range("D35").Value = '0.125'
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.38'
$ws.Range("E37").Value = '  -3.04%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0350'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.60'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +27.63%  '
$ws.Range("D40").Value = '0.251'
$ws.Range("E40").Value = '  +26.48%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.50'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").Value = '2.18'
$ws.Range("E42").Value = '  +4.22%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = '62.13'
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = '5.43'
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("E45").Value = '  +4.56%  '
$ws.Range("D46").Value = '8.42'
$ws.Range("E46").Value = '  +1.55%  '
$ws.Range("D47").Value = '100.65'
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("E48").Value = '  +3.99%  '
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").Value = '2.68'
$ws.Range("E51").Value = '  +1.09%  '
